$d = $word.ActiveDocument

# Remove the "Line 2: Pressed either the Start button, or another button
# that is linked right after the jump to perform a combo." sentence from
# the Alternative Flow of Events row, leaving the leading space run intact.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$result = $find.Execute(
    "Line 2: Pressed either the Start button, or another button that is linked right after the jump to perform a combo.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2
)
